# Applies the cryptos.xlsx price/volume update described in the commit diff.
# Numeric-looking text values (e.g. "210.10") are written with a leading
# apostrophe to force Excel to keep them as text (matching the workbook's
# original inline-string cells), then the cell style is reset to "Normal"
# so no stray text/quote-prefix formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.323.82'
$ws.Range("E2").Value = '  -1.12%  '

# Row 3
$ws.Range("D3").Value = '1.586.90'
$ws.Range("E3").Value = '  -0.56%  '

# Row 4
$ws.Range("E4").Value = '  -0.41%  '

# Row 5
$ws.Range("D5").Value = '''210.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.19%  '

# Row 6
$ws.Range("E6").Value = '  -1.08%  '

# Row 7
$ws.Range("E7").Value = '  -0.42%  '

# Row 8
$ws.Range("E8").Value = '  -0.65%  '

# Row 9
$ws.Range("E9").Value = '  -0.46%  '

# Row 10
$ws.Range("D10").Value = '''19.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.26%  '

# Row 11
$ws.Range("D11").Value = '''0.0844'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.11%  '

# Row 12
$ws.Range("D12").Value = '1.811.61'
$ws.Range("E12").Value = '  -0.48%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.591.17'
$ws.Range("E13").Value = '  +0.53%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.81%  '

# Row 15
$ws.Range("E15").Value = '  -0.69%  '

# Row 16
$ws.Range("D16").Value = '''64.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.45%  '

# Row 17
$ws.Range("D17").Value = '26.328.02'
$ws.Range("E17").Value = '  -1.01%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0726'
$ws.Range("E18").Value = '  -1.56%  '

# Row 19
$ws.Range("D19").Value = '''7.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.09%  '

# Row 20
$ws.Range("D20").Value = '''210.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.40%  '

# Row 21
$ws.Range("E21").Value = '  -0.41%  '

# Row 22
$ws.Range("D22").Value = '''4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '

# Row 23
$ws.Range("E23").Value = '  -3.32%  '

# Row 24
$ws.Range("E24").Value = '  -0.09%  '

# Row 25
$ws.Range("D25").Value = '''144.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.79%  '

# Row 27
$ws.Range("E27").Value = '  -1.10%  '

# Row 28
$ws.Range("E28").Value = '  -0.41%  '

# Row 29
$ws.Range("D29").Value = '''15.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.27%  '

# Row 30
$ws.Range("E30").Value = '  -0.09%  '

# Row 31
$ws.Range("E31").Value = '  -0.30%  '

# Row 32
$ws.Range("E32").Value = '  -0.91%  '

# Row 33
$ws.Range("E33").Value = '  +1.25%  '

# Row 34
$ws.Range("D34").Value = '1.301.12'
$ws.Range("E34").Value = '  +1.75%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.66%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.611'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.64%  '

# Row 37
$ws.Range("E37").Value = '  -0.71%  '

# Row 38
$ws.Range("E38").Value = '  +0.39%  '

# Row 39
$ws.Range("E39").Value = '  -12.72%  '

# Row 40
$ws.Range("D40").Value = '''0.807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.80%  '

# Row 41
$ws.Range("E41").Value = '  -0.39%  '

# Row 42
$ws.Range("E42").Value = '  +3.40%  '

# Row 43
$ws.Range("E43").Value = '  -0.39%  '

# Row 44
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '''2.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.24%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''62.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.15%  '

# Row 46
$ws.Range("D46").Value = '1.723.31'
$ws.Range("E46").Value = '  -0.53%  '

# Row 47
$ws.Range("D47").Value = '''87.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.90%  '

# Row 48
$ws.Range("D48").Value = '''1.48'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.49%  '

# Row 49
$ws.Range("E49").Value = '  -1.51%  '

# Row 50
$ws.Range("E50").Value = '  -4.37%  '

# Row 51
$ws.Range("E51").Value = '  -0.34%  '
